$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "MCT-2A-Circuitos elétricos 2"
$ws.Range("F2").Value = "-"

$ws.Range("C3").Value = "-"

$ws.Range("C4").Value = "[-, -, -, 'MCT-2A-Programação de Computadores']"

$ws.Range("C6").Value = "[-, -, -, 'MCT-2A-Programação de Computadores']"

$ws.Range("C7").Value = "[-, -, -, 'MCT-2A-Programação de Computadores']"

$ws.Range("B8").Value = "MCT-2A-Circuitos elétricos 2"
$ws.Range("C8").Value = "[-, -, -, 'MCT-2A-Programação de Computadores']"
$ws.Range("F8").Value = "-"

$ws.Range("B18").Value = "-"
$ws.Range("E18").Value = "-"

$ws.Range("B19").Value = "-"
$ws.Range("D19").Value = "-"
$ws.Range("E19").Value = "-"

$ws.Range("B20").Value = "-"
$ws.Range("D20").Value = "-"

$ws.Range("C21").Value = "-"
$ws.Range("D21").Value = "-"
$ws.Range("E21").Value = "-"
